$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F (dSF)
$updates = @{
    2  = -4
    3  = -3
    4  = -4
    5  = 2
    6  = -4
    7  = -1
    11 = -1
    12 = -1
    13 = -7
    14 = 1
    16 = -5
    19 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
